# Adds the new intro paragraphs + "15 Animation for Doors" Heading1
# paragraph + a trailing blank paragraph right after the "Write Up"
# title paragraph, matching the target diff exactly. We build the new
# content as a Word-flavoured-OOXML fragment (so we can reproduce the
# <w:proofErr .../> grammar-check markers around "don't" and the
# xml:space="preserve" runs) and drop it in with Range.InsertXML,
# which is the COM-exposed way to inject raw WordOpenXML.

$d = $word.ActiveDocument

# Collapsed range sitting right at the end of the "Write Up" paragraph
# (i.e. just before its paragraph mark) -- inserting OOXML paragraphs
# here splits cleanly into new <w:p> elements without disturbing the
# existing empty paragraph / Heading1 paragraph that already follow it.
$titlePara = $d.Paragraphs(1)
$insertionPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

$newContentXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>In this tutorial, we will be looking at writing the code to give our doors some special effects. Yes, we will be animating those effects right in Game Maker’s code panel and giving them a bit of razzle-dazzle.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">So, if you would like to learn a bit more on how to go about doing this, then why </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>don’t</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> you join us for our brand-new article entitled:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
            </w:pPr>
            <w:r>
              <w:t>15 Animation for Doors</w:t>
            </w:r>
          </w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$insertionPoint.InsertXML($newContentXml)
